$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.418.58'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '2.635.60'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.38'
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.35'
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.558'
$ws.Range("E8").Value = '  +3.07%  '
$ws.Range("D9").Value = '2.636.29'
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("E10").Value = '  +3.96%  '
$ws.Range("E11").Value = '  +0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.21'
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("E13").Value = '  -1.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.80'
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").Value = '3.116.02'
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").Value = '67.363.69'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '2.637.09'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '364.14'
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.57'
$ws.Range("E21").Value = '  -3.96%  '
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.13'
$ws.Range("E23").Value = '  +3.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '66.35'
$ws.Range("E26").Value = '  -7.23%  '
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '579.35'
$ws.Range("E30").Value = '  -7.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.41'
$ws.Range("E31").Value = '  -4.11%  '
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.129'
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("E36").Value = '  -2.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.93'
$ws.Range("E37").Value = '  -1.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.03'
$ws.Range("E38").Value = '  +2.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.46'
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.30'
$ws.Range("E41").Value = '  -3.84%  '
$ws.Range("E42").Value = '  -1.31%  '
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.20'
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '156.09'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.74'
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.626'
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.74'
$ws.Range("E51").Value = '  -1.59%  '
